$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header B1 text, remove C1/D1 (nb_dos, nb_file) columns entirely
$ws.Range("B1").Value = "Nombre de fichiers dans le dossier"

# Delete columns C and D (which held nb_dos / nb_file)
$ws.Range("C1:D1").EntireColumn.Delete()

# Set column B width to fit content (bestFit-like behavior)
$ws.Columns.Item(2).AutoFit()

# Adjust window size/position to match target
$win = $excel.ActiveWindow
$win.Left = 3180
$win.Top = 1440
$win.Width = 10965
$win.Height = 9990
